$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1717.25
$ws.Range("J18").Value = 2225
$ws.Range("L18").Value = 2225
$ws.Range("N18").Value = -2793

# Row 41
$ws.Range("H41").Value = 1987.238
$ws.Range("I41").Value = 1131.6471
$ws.Range("K41").Value = 1131.6471
$ws.Range("M41").Value = -691.6470999999999

# Row 43
$ws.Range("H43").Value = 5840.4
$ws.Range("I43").Value = 2500
$ws.Range("J43").Value = 8067.3335
$ws.Range("K43").Value = 2500
$ws.Range("L43").Value = 8067.3335
$ws.Range("M43").Value = -2431
$ws.Range("N43").Value = -8205.333500000001

# Row 47
$ws.Range("H47").Value = 300000
$ws.Range("I47").Value = 300000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 300000
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -299028

# Row 70
$ws.Range("H70").Value = 4953.778
$ws.Range("I70").Value = 4346.6665
$ws.Range("J70").Value = 5560.8887
$ws.Range("K70").Value = 13039.9995
$ws.Range("L70").Value = 16682.6661
$ws.Range("M70").Value = -12769.9995
$ws.Range("N70").Value = -17222.6661

# Row 73
$ws.Range("H73").Value = 4953.778
$ws.Range("I73").Value = 4346.6665
$ws.Range("J73").Value = 5560.8887
$ws.Range("K73").Value = 13039.9995
$ws.Range("L73").Value = 16682.6661
$ws.Range("M73").Value = -12103.9995
$ws.Range("N73").Value = -18554.6661

# Row 76
$ws.Range("H76").Value = 13334.667
$ws.Range("I76").Value = 10500
$ws.Range("K76").Value = 10500
$ws.Range("M76").Value = -10185

# Row 79
$ws.Range("H79").Value = 13334.667
$ws.Range("I79").Value = 10500
$ws.Range("K79").Value = 10500
$ws.Range("M79").Value = -9408

# Row 80
$ws.Range("H80").Value = 2561.2144
$ws.Range("I80").Value = 432.5
$ws.Range("J80").Value = 3412.7
$ws.Range("K80").Value = 1297.5
$ws.Range("L80").Value = 10238.1
$ws.Range("M80").Value = -299.5
$ws.Range("N80").Value = -12234.1

# Row 83
$ws.Range("H83").Value = 2561.2144
$ws.Range("I83").Value = 432.5
$ws.Range("J83").Value = 3412.7
$ws.Range("K83").Value = 3892.5
$ws.Range("L83").Value = 30714.3
$ws.Range("M83").Value = 1099.5
$ws.Range("N83").Value = -40698.3

# Row 86
$ws.Range("H86").Value = 6988.75
$ws.Range("I86").Value = 7350
$ws.Range("J86").Value = 6627.5
$ws.Range("K86").Value = 7350
$ws.Range("L86").Value = 6627.5
$ws.Range("M86").Value = -6227
$ws.Range("N86").Value = -8873.5

# Row 89
$ws.Range("H89").Value = 6988.75
$ws.Range("I89").Value = 7350
$ws.Range("J89").Value = 6627.5
$ws.Range("K89").Value = 36750
$ws.Range("L89").Value = 33137.5
$ws.Range("M89").Value = -31134
$ws.Range("N89").Value = -44369.5

# Row 113
$ws.Range("H113").Value = 19006
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 19006
$ws.Range("K113").Value = 0
$ws.Range("L113").ClearContents()
$ws.Range("M113").Value = 19006
$ws.Range("N113").Value = -25514

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = 0

# Row 137
$ws.Range("H137").Value = 2682
$ws.Range("I137").Value = 687.1429000000001
$ws.Range("K137").Value = 2061.4287
$ws.Range("M137").Value = 488.5712999999996

# Row 138
$ws.Range("H138").Value = 2712.2295
$ws.Range("I138").Value = 1682.742
$ws.Range("J138").Value = 3776.0334
$ws.Range("K138").Value = 5048.226
$ws.Range("L138").Value = 11328.1002
$ws.Range("M138").Value = 91.77400000000034
$ws.Range("N138").Value = -21608.1002


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3553.3948
$ws.Range("I32").Value = 2843.0286
$ws.Range("K32").Value = 2843.0286
$ws.Range("M32").Value = -2556.0286

# Row 122
$ws.Range("H122").Value = 5193.625
$ws.Range("J122").Value = 6624.75
$ws.Range("L122").Value = 19874.25
$ws.Range("N122").Value = -24774.25

# Row 132
$ws.Range("H132").Value = 7966.737
$ws.Range("J132").Value = 9306.846
$ws.Range("L132").Value = 27920.538
$ws.Range("N132").Value = -32980.538


$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

# Row 64
$ws.Range("H64").Value = 3056.8
$ws.Range("I64").Value = 2103
$ws.Range("J64").Value = 3692.6667
$ws.Range("K64").Value = 2103
$ws.Range("L64").Value = 3692.6667
$ws.Range("M64").Value = -1878
$ws.Range("N64").Value = -4142.6667

# Row 67
$ws.Range("H67").Value = 3056.8
$ws.Range("I67").Value = 2103
$ws.Range("J67").Value = 3692.6667
$ws.Range("K67").Value = 2103
$ws.Range("L67").Value = 3692.6667
$ws.Range("M67").Value = -1323
$ws.Range("N67").Value = -5252.6667

# Row 105
$ws.Range("H105").Value = 20979.785
$ws.Range("J105").Value = 9683.5
$ws.Range("L105").Value = 9683.5
$ws.Range("N105").Value = -13177.5


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5460.2266
$ws.Range("I31").Value = 3204.6924
$ws.Range("J31").Value = 11743.5
$ws.Range("K31").Value = 3204.6924
$ws.Range("L31").Value = 11743.5
$ws.Range("M31").Value = -2909.6924
$ws.Range("N31").Value = -12333.5

# Row 34
$ws.Range("H34").Value = 5460.2266
$ws.Range("I34").Value = 3204.6924
$ws.Range("J34").Value = 11743.5
$ws.Range("K34").Value = 3204.6924
$ws.Range("L34").Value = 11743.5
$ws.Range("M34").Value = -3002.6924
$ws.Range("N34").Value = -12147.5

# Row 107
$ws.Range("H107").Value = 1195.3334
$ws.Range("I107").Value = 1145.5834
$ws.Range("K107").Value = 1145.5834
$ws.Range("M107").Value = 774.4166

# Row 122
$ws.Range("H122").Value = 9721.888999999999
$ws.Range("I122").Value = 5333
$ws.Range("K122").Value = 15999
$ws.Range("M122").Value = -13549

# Row 132
$ws.Range("H132").Value = 3828.2856
$ws.Range("I132").Value = 3477.0227
$ws.Range("J132").Value = 5116.25
$ws.Range("K132").Value = 10431.0681
$ws.Range("L132").Value = 15348.75
$ws.Range("M132").Value = -7901.0681
$ws.Range("N132").Value = -20408.75

# Row 134
$ws.Range("H134").Value = 2706.9285
$ws.Range("I134").Value = 1716.7273
$ws.Range("J134").Value = 6337.6665
$ws.Range("K134").Value = 5150.1819
$ws.Range("L134").Value = 19012.9995
$ws.Range("M134").Value = -2615.1819
$ws.Range("N134").Value = -24082.9995


$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 3600
$ws.Range("I41").Value = 3600
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3600
$ws.Range("L41").ClearContents()
$ws.Range("M41").Value = -3245
$ws.Range("N41").Value = 0

# Row 80
$ws.Range("H80").Value = 10499.75
$ws.Range("I80").Value = 8249.5
$ws.Range("J80").Value = 12750
$ws.Range("K80").Value = 8249.5
$ws.Range("L80").Value = 12750
$ws.Range("M80").Value = -7251.5
$ws.Range("N80").Value = -14746

# Row 83
$ws.Range("H83").Value = 10499.75
$ws.Range("I83").Value = 8249.5
$ws.Range("J83").Value = 12750
$ws.Range("K83").Value = 41247.5
$ws.Range("L83").Value = 63750
$ws.Range("M83").Value = -36255.5
$ws.Range("N83").Value = -73734

# Row 122
$ws.Range("H122").Value = 7050.875
$ws.Range("I122").Value = 6399.8335
$ws.Range("J122").Value = 9004
$ws.Range("K122").Value = 19199.5005
$ws.Range("L122").Value = 27012
$ws.Range("M122").Value = -16749.5005
$ws.Range("N122").Value = -31912

# Row 132
$ws.Range("H132").Value = 48807.043
$ws.Range("I132").Value = 51697.953
$ws.Range("K132").Value = 155093.859
$ws.Range("M132").Value = -152563.859


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2183.4546
$ws.Range("I16").Value = 2259.4
$ws.Range("J16").Value = 1424
$ws.Range("K16").Value = 2259.4
$ws.Range("L16").Value = 1424
$ws.Range("M16").Value = -2089.4
$ws.Range("N16").Value = -1764

# Row 42
$ws.Range("H42").Value = 40012.5
$ws.Range("I42").Value = 40025
$ws.Range("J42").Value = 40000
$ws.Range("K42").Value = 40025
$ws.Range("L42").Value = 40000
$ws.Range("M42").Value = -39462
$ws.Range("N42").Value = -41126

# Row 46
$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

# Row 49
$ws.Range("H49").Value = 40012.5
$ws.Range("I49").Value = 40025
$ws.Range("J49").Value = 40000
$ws.Range("K49").Value = 40025
$ws.Range("L49").Value = 40000
$ws.Range("M49").Value = -39878
$ws.Range("N49").Value = -40294

# Row 122
$ws.Range("H122").Value = 8483.777
$ws.Range("I122").Value = 6887.25
$ws.Range("J122").Value = 9761
$ws.Range("K122").Value = 20661.75
$ws.Range("L122").Value = 29283
$ws.Range("M122").Value = -18211.75
$ws.Range("N122").Value = -34183


$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").ClearContents()
$ws.Range("N38").Value = 0

# Row 122
$ws.Range("H122").Value = 12743.723
$ws.Range("I122").Value = 4375.4
$ws.Range("J122").Value = 15962.308
$ws.Range("K122").Value = 13126.2
$ws.Range("L122").Value = 47886.924
$ws.Range("M122").Value = -10676.2
$ws.Range("N122").Value = -52786.924

# Row 132
$ws.Range("H132").Value = 3767.4482
$ws.Range("I132").Value = 3653.7778
$ws.Range("K132").Value = 10961.3334
$ws.Range("M132").Value = -8431.3334

# Row 135
$ws.Range("H135").Value = 55817.363
$ws.Range("J135").Value = 55817.363
$ws.Range("L135").Value = 55817.363
$ws.Range("N135").Value = -65957.363

# Row 136
$ws.Range("H136").Value = 3574.3276
$ws.Range("I136").Value = 2345.451
$ws.Range("K136").Value = 7036.353
$ws.Range("M136").Value = -4486.353

